$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect, make the edits, then re-protect so the
# workbook's protection state is restored afterwards.
$ws.Unprotect()

# Update the confidential notice date from 2021-05-03 to 2021-05-04
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-04 for illustrative purposes only and are subject to change."

# Update Weight / Percent Change figures
$ws.Range("D2").Value = 0.8458231879297664
$ws.Range("E2").Value = -0.01481575281752556

$ws.Range("D3").Value = 0.1541768120702336
$ws.Range("E3").Value = -0.01130048165987396

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = -0.01427377951687647

$ws.Protect()
